# Add two new rows (159, 160) of parsed data to each of the four sheets
# (MID_LFT_#1, MID_LFT_#2, MID_PLT_#1, MID_PLT_#2), continuing the daily
# log pattern that ends at row 158.

$wb = $excel.ActiveWorkbook

# Per-sheet new-row data: time(A), B, C, D(row159), D(row160), E, F, G, H(row159), H(row160), I
$sheetsData = @(
    @{
        Index = 1
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D159 = "0x00,0xC4"
        D160 = "0x00,0xC0"
        E = "0x07"
        F = 400
        G = "5.68631262647113e+23"
        H159 = 204
        H160 = 204
        I159 = 7
        I160 = 7
    },
    @{
        Index = 2
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D159 = "0x00,0xE0"
        D160 = "0x00,0xDC"
        E = "0x19"
        F = 380
        G = "5.68432987514711e+23"
        H159 = 224
        H160 = 220
        I159 = 25
        I160 = 25
    },
    @{
        Index = 3
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D159 = "0x00,0x51"
        D160 = "0x00,0x51"
        E = "0x15"
        F = 110
        G = "5.68631262647113e+23"
        H159 = 81
        H160 = 81
        I159 = 15
        I160 = 15
    },
    @{
        Index = 4
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D159 = "0x00,0x6A"
        D160 = "0x00,0x69"
        E = "0x9"
        F = 130
        G = "5.68631262647113e+23"
        H159 = 106
        H160 = 105
        I159 = 9
        I160 = 9
    }
)

$dateRow159 = 45945.46252314815
$dateRow160 = 45946.46368055556
$dateFormat = "YYYY-MM-DD HH:MM:SS"

foreach ($sd in $sheetsData) {
    $ws = $wb.Worksheets.Item($sd.Index)

    # Row 159
    $ws.Cells.Item(159, 1).Value = $dateRow159
    $ws.Cells.Item(159, 1).NumberFormat = $dateFormat
    $ws.Cells.Item(159, 2).Value = $sd.B
    $ws.Cells.Item(159, 3).Value = $sd.C
    $ws.Cells.Item(159, 4).Value = $sd.D159
    $ws.Cells.Item(159, 5).Value = $sd.E
    $ws.Cells.Item(159, 6).Value = $sd.F
    $ws.Cells.Item(159, 7).Value = [double]$sd.G
    $ws.Cells.Item(159, 8).Value = $sd.H159
    $ws.Cells.Item(159, 9).Value = $sd.I159

    # Row 160
    $ws.Cells.Item(160, 1).Value = $dateRow160
    $ws.Cells.Item(160, 1).NumberFormat = $dateFormat
    $ws.Cells.Item(160, 2).Value = $sd.B
    $ws.Cells.Item(160, 3).Value = $sd.C
    $ws.Cells.Item(160, 4).Value = $sd.D160
    $ws.Cells.Item(160, 5).Value = $sd.E
    $ws.Cells.Item(160, 6).Value = $sd.F
    $ws.Cells.Item(160, 7).Value = [double]$sd.G
    $ws.Cells.Item(160, 8).Value = $sd.H160
    $ws.Cells.Item(160, 9).Value = $sd.I160
}
